# validate with test data
#
# The RESULT rows (the 4th row of every Bayes/Idtax/LCA/RESULT quadruplet)
# were missing their taxonomy-rank values in columns D:K. This fills those
# consensus values back in so each RESULT row matches the taxonomy chosen
# for that ASV, exactly as produced by the pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 66 (RESULT, sv19589 group)
$ws.Range("D66").Value = "Eukaryota"
$ws.Range("E66").Value = "Alveolata"
$ws.Range("F66").Value = "Dinoflagellata"
$ws.Range("G66").Value = "Syndiniales"
$ws.Range("H66").Value = "NA"

# Row 86 (RESULT, sv17874 group)
$ws.Range("D86").Value = "Eukaryota"
$ws.Range("E86").Value = "Alveolata"
$ws.Range("F86").Value = "Dinoflagellata"
$ws.Range("G86").Value = "Dinophyceae"
$ws.Range("H86").Value = "NA"
$ws.Range("I86").Value = "NA"
$ws.Range("J86").Value = "NA"

# Row 91 (RESULT)
$ws.Range("D91").Value = "NA"

# Row 121 (RESULT)
$ws.Range("D121").Value = "Eukaryota"
$ws.Range("E121").Value = "NA"
$ws.Range("F121").Value = "NA"
$ws.Range("G121").Value = "NA"
$ws.Range("H121").Value = "NA"
$ws.Range("I121").Value = "NA"
$ws.Range("J121").Value = "NA"

# Row 131 (RESULT, Arthropoda group)
$ws.Range("D131").Value = "Eukaryota"
$ws.Range("E131").Value = "Opisthokonta"
$ws.Range("F131").Value = "Metazoa"
$ws.Range("G131").Value = "Arthropoda"
$ws.Range("H131").Value = "NA"
$ws.Range("I131").Value = "NA"

# Row 136 (RESULT)
$ws.Range("D136").Value = "Eukaryota"
$ws.Range("E136").Value = "Opisthokonta"
$ws.Range("F136").Value = "Metazoa"
$ws.Range("G136").Value = "NA"

# Row 141 (RESULT, Arthropoda group)
$ws.Range("D141").Value = "Eukaryota"
$ws.Range("E141").Value = "Opisthokonta"
$ws.Range("F141").Value = "Metazoa"
$ws.Range("G141").Value = "Arthropoda"
$ws.Range("H141").Value = "NA"
$ws.Range("I141").Value = "NA"

# Row 146 (RESULT)
$ws.Range("D146").Value = "Eukaryota"
$ws.Range("E146").Value = "NA"
$ws.Range("F146").Value = "NA"

# Row 151 (RESULT, Ascomycota/Pezizomycotina group)
$ws.Range("D151").Value = "Eukaryota"
$ws.Range("E151").Value = "Opisthokonta"
$ws.Range("F151").Value = "Fungi"
$ws.Range("G151").Value = "Ascomycota"
$ws.Range("H151").Value = "Pezizomycotina"
$ws.Range("I151").Value = "NA"

# Row 161 (RESULT)
$ws.Range("D161").Value = "Eukaryota"
$ws.Range("E161").Value = "NA"
$ws.Range("F161").Value = "NA"
$ws.Range("G161").Value = "NA"
$ws.Range("H161").Value = "NA"

# Row 166 (RESULT, Dothideomycetes group)
$ws.Range("D166").Value = "Eukaryota"
$ws.Range("E166").Value = "Opisthokonta"
$ws.Range("F166").Value = "Fungi"
$ws.Range("G166").Value = "Ascomycota"
$ws.Range("H166").Value = "Pezizomycotina"
$ws.Range("I166").Value = "Dothideomycetes"
$ws.Range("J166").Value = "NA"

# Row 171 (RESULT)
$ws.Range("D171").Value = "Eukaryota"
$ws.Range("E171").Value = "NA"
$ws.Range("F171").Value = "NA"
$ws.Range("G171").Value = "NA"
$ws.Range("H171").Value = "NA"
$ws.Range("I171").Value = "NA"
$ws.Range("J171").Value = "NA"

# Row 186 (RESULT, Dothideomycetes group)
$ws.Range("D186").Value = "Eukaryota"
$ws.Range("E186").Value = "Opisthokonta"
$ws.Range("F186").Value = "Fungi"
$ws.Range("G186").Value = "Ascomycota"
$ws.Range("H186").Value = "Pezizomycotina"
$ws.Range("I186").Value = "Dothideomycetes"
$ws.Range("J186").Value = "NA"
$ws.Range("K186").Value = "NA"

# Row 191 (RESULT, Dothideomycetes group)
$ws.Range("D191").Value = "Eukaryota"
$ws.Range("E191").Value = "Opisthokonta"
$ws.Range("F191").Value = "Fungi"
$ws.Range("G191").Value = "Ascomycota"
$ws.Range("H191").Value = "Pezizomycotina"
$ws.Range("I191").Value = "Dothideomycetes"
$ws.Range("J191").Value = "NA"

# Update the sheet view: drop the frozen/scrolled top-left cell, change the
# zoom level, and move the active selection, matching the sheetView in the
# diff (zoomScale 138 -> 116, selection J61 -> F195, no more topLeftCell).
$excel.ActiveWindow.Zoom = 116
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F195").Select()
